$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Friday (2/26) time tracking -- enter time for week 2 (row 4): Tuesday = 15, Friday = 55
$ws.Range("D4").Value = 15
$ws.Range("G4").Value = 55

# Move the active selection to G5, matching where the user clicked next
$ws.Range("G5").Select()
